# Append the newly scraped price row (2026-02-07) to the tracking sheet.
# The existing sheet stores every Date/Price/Discount/Incredible value as
# plain shared-string text (not as native numbers/dates), so the new cells
# must be written the same way: as literal text, with no number formatting
# and no style applied to the cell (matching the rest of the sheet, which
# uses the default style only).
#
# Directly assigning a numeric- or date-looking string to .Value makes
# Excel auto-convert it to a real number/date, and prefixing with an
# apostrophe to force text marks the cell with a "quote prefix" style -
# neither matches the target. Instead we build the literal text with a
# formula (="...") in a scratch cell, copy it, and paste-special just the
# values into the destination; this keeps the result as plain shared-string
# text without attaching any style to the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newRow = 38
$scratch = $ws.Cells.Item(1, 20)   # far-away unused cell used as scratch space

function Set-TextValue {
    param($cell, [string]$text)

    $scratch.Formula = '="' + $text + '"'
    $scratch.Copy()
    $cell.PasteSpecial(-4163)       # xlPasteValues
}

Set-TextValue $ws.Cells.Item($newRow, 1) "2026-02-07"
Set-TextValue $ws.Cells.Item($newRow, 2) "336400"
Set-TextValue $ws.Cells.Item($newRow, 3) "47"
Set-TextValue $ws.Cells.Item($newRow, 4) "0"

$scratch.Clear()
